# Connectathon marking.xlsx -- add a "scorecard" sheet, used to record
# participant / scenario results, and retarget the stray "Clients" label
# on "summary by scenario" (it now belongs to the new sheet).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Add the new "scorecard" worksheet as the last tab.
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "scorecard"

# Column widths (character units; close to the authored sheet).
$ws.Columns.Item(2).ColumnWidth = 10.83
$ws.Columns.Item(3).ColumnWidth = 7.5
$ws.Columns.Item(4).ColumnWidth = 4.66
$ws.Columns.Item(5).ColumnWidth = 9.83
$ws.Columns.Item(6).ColumnWidth = 6.66
$ws.Columns.Item(7).ColumnWidth = 4.66
$ws.Columns.Item(8).ColumnWidth = 7
$ws.Columns.Item(9).ColumnWidth = 9.83
$ws.Columns.Item(10).ColumnWidth = 6.83
$ws.Columns.Item(11).ColumnWidth = 6.83
$ws.Columns.Item(12).ColumnWidth = 5.83
$ws.Columns.Item(13).ColumnWidth = 9.83
$ws.Columns.Item(14).ColumnWidth = 7.83
$ws.Columns.Item(15).ColumnWidth = 5.66
$ws.Columns.Item(16).ColumnWidth = 6.83
$ws.Columns.Item(17).ColumnWidth = 9.83
$ws.Columns.Item(19).ColumnWidth = 6.5
$ws.Columns.Item(20).ColumnWidth = 7.16

# ---------------------------------------------------------------------
# 2. The label that used to live on "summary by scenario" now belongs
#    on the new sheet, so retarget that cell to a real scenario/
#    participant label first (this is the first *new* string used).
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("summary by scenario")
$summary.Range("C5").Value = "participant 1"

# ---------------------------------------------------------------------
# 3. Fill in the participant list (left column) next ...
# ---------------------------------------------------------------------
$ws.Range("B5").Value = "Clients"
$ws.Range("B5").Font.Bold = $true

$ws.Range("B8").Value = "David Browser"
$ws.Range("B8").Font.Bold = $true
$ws.Range("B9").Value = "David Android"
$ws.Range("B9").Font.Bold = $true
$ws.Range("B10").Value = "Claude"
$ws.Range("B10").Font.Bold = $true
$ws.Range("B11").Value = "Rolim"
$ws.Range("B11").Font.Bold = $true
$ws.Range("B12").Value = "Healthfile"
$ws.Range("B12").Font.Bold = $true
$ws.Range("B13").Value = "Eliot"
$ws.Range("B13").Font.Bold = $true

# ---------------------------------------------------------------------
# 3. ... then the four scenario headers ...
# ---------------------------------------------------------------------
$ws.Range("E4").Value = "Scenario1"
$ws.Range("E4").Font.Bold = $true
$ws.Range("I4").Value = "Scenario2"
$ws.Range("I4").Font.Bold = $true
$ws.Range("M4").Value = "Scenario3"
$ws.Range("M4").Font.Bold = $true
$ws.Range("Q4").Value = "Scenario4"
$ws.Range("Q4").Font.Bold = $true

# ---------------------------------------------------------------------
# 4. ... then the big title ...
# ---------------------------------------------------------------------
$ws.Range("A2").Value = "Result scorecard"
$ws.Range("A2").Font.Bold = $true
$ws.Range("A2").Font.Size = 18
$ws.Rows.Item(2).RowHeight = 23

# ---------------------------------------------------------------------
# 5. ... then the footnote, the "Notes" column header and its sample note.
# ---------------------------------------------------------------------
$ws.Range("A18").Value = "Place the representation (json / xml) in the cell"

$ws.Range("U4").Value = "Notes"
$ws.Range("U4").Font.Bold = $true

$ws.Range("U8").Value = "Claim extensions in problem processing"

# ---------------------------------------------------------------------
# 6. Server-name sub-headers, repeated under each scenario block.
# ---------------------------------------------------------------------
foreach ($col in @("E", "I", "M", "Q")) {
    $ws.Range($col + "5").Value = "Grahame"
}
foreach ($col in @("F", "J", "N", "R")) {
    $ws.Range($col + "5").Value = "Ewout"
}
foreach ($col in @("G", "K", "O", "S")) {
    $ws.Range($col + "5").Value = "Rik"
}
foreach ($col in @("H", "L", "P", "T")) {
    $ws.Range($col + "5").Value = "Brett"
}

# Row 8 / row 9 placeholders ("json") under the first two participants.
foreach ($col in @("E", "F", "I", "J", "M", "N", "Q", "R", "T")) {
    $ws.Range($col + "8").Value = "json"
}
foreach ($col in @("I", "J", "Q", "R", "T")) {
    $ws.Range($col + "9").Value = "json"
}

# ---------------------------------------------------------------------
# 7. Scenario legend at the bottom of the sheet.
# ---------------------------------------------------------------------
$ws.Range("A21").Value = "Scenarios"
$ws.Range("A21").Font.Bold = $true

$ws.Range("A23").Value = 1
$ws.Range("B23").Value = "register new patient"
$ws.Range("A24").Value = 2
$ws.Range("B24").Value = "update patient"
$ws.Range("A25").Value = 3
$ws.Range("B25").Value = "retrieve history"
$ws.Range("A26").Value = 4
$ws.Range("B26").Value = "patient search"

# ---------------------------------------------------------------------
# 8. Page setup to match the workbook's other sheets.
# ---------------------------------------------------------------------
$ps = $ws.PageSetup
$ps.LeftMargin = 54
$ps.RightMargin = 54
$ps.TopMargin = 72
$ps.BottomMargin = 72
$ps.HeaderMargin = 36
$ps.FooterMargin = 36
$ps.PaperSize = 9
$ps.Orientation = 1

# ---------------------------------------------------------------------
# 9. Selection state on "summary by scenario".
# ---------------------------------------------------------------------
$summary.Activate()
$summary.Range("C5").Select()

# ---------------------------------------------------------------------
# 10. Selection / scroll state, and make "scorecard" the active tab.
# ---------------------------------------------------------------------
$registration = $wb.Worksheets.Item("registration")
$registration.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 2
$win.ScrollColumn = 1
$registration.Range("E15").Select()

$ws.Activate()
$ws.Range("K24").Select()
